$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45..161 down to 46..162
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new price record
$ws.Range("A45").Value = 6
$ws.Range("B45").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44495
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 100112026
$ws.Range("G45").Value = "Haba"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 4000
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = 4271
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 171
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
